$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the "Förändrad" date serial value by 1 day (46061 -> 46062)
# for rows 2 through 14 in column C.
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Range("C$row")
    $cell.Value2 = $cell.Value2 + 1
}
